$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 19:43"

# Country name swaps (ranking shuffled after the data refresh) and updated
# COVID-19 statistics (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) for the affected rows.
$updates = @(
    ,("B4", 2194477)
    ,("C4", 11527)
    ,("D4", 892093)
    ,("E4", 1183682)
    ,("G4", 419)
    ,("H4", 118702)
    ,("B5", 904734)
    ,("C5", 13178)
    ,("E5", 395303)
    ,("G5", 539)
    ,("H5", 44657)
    ,("B7", 352815)
    ,("C7", 9789)
    ,("D7", 187356)
    ,("E7", 153577)
    ,("G7", 1967)
    ,("H7", 11882)
    ,("B13", 188331)
    ,("C13", 287)
    ,("E13", 6332)
    ,("G13", 14)
    ,("H13", 8899)
    ,("B15", 181298)
    ,("C15", 1467)
    ,("D15", 153379)
    ,("E15", 23077)
    ,("G15", 17)
    ,("H15", 4842)
    ,("B39", 31154)
    ,("C39", 23)
    ,("E39", 300)
    ,("G39", 15)
    ,("H39", 1954)
    ,("B43", 25334)
    ,("C43", 13)
    ,("E43", 927)
    ,("G43", 3)
    ,("H43", 1709)
    ,("B68", 8931)
    ,("C68", 46)
    ,("D68", 7937)
    ,("E68", 782)
    ,("D76", 4096)
    ,("E76", 1213)
    ,("B81", 4539)
    ,("C81", 38)
    ,("D81", 3324)
    ,("E81", 1172)
    ,("B85", 4075)
    ,("C85", 3)
    ,("D85", 3933)
    ,("E85", 32)
    ,("B95", 2658)
    ,("C95", 16)
    ,("D95", 649)
    ,("E95", 1921)
    ,("A97", "Republica de Africa Central")
    ,("B97", 2410)
    ,("C97", 188)
    ,("D97", 396)
    ,("E97", 2000)
    ,("G97", 7)
    ,("H97", 14)
    ,("A98", "Mayotte")
    ,("B98", 2333)
    ,("C98", 23)
    ,("D98", 2058)
    ,("E98", 246)
    ,("G98", 0)
    ,("H98", 29)
    ,("B101", 2094)
    ,("C101", 29)
    ,("D101", 1670)
    ,("E101", 416)
    ,("B103", 1915)
    ,("C103", 10)
    ,("E103", 533)
    ,("B122", 1225)
    ,("C122", 49)
    ,("D122", 686)
    ,("E122", 488)
    ,("A148", "Suazilandia")
    ,("B148", 520)
    ,("C148", 14)
    ,("D148", 259)
    ,("E148", 257)
    ,("H148", 4)
    ,("A149", "Estado de Palestina")
    ,("B149", 511)
    ,("C149", 6)
    ,("D149", 415)
    ,("E149", 93)
    ,("H149", 3)
    ,("A150", "Tanzania")
    ,("B150", 509)
    ,("D150", 183)
    ,("E150", 305)
    ,("H150", 21)
    ,("A151", "Liberia")
    ,("B151", 509)
    ,("C151", 11)
    ,("D151", 222)
    ,("E151", 254)
    ,("H151", 33)
    ,("D160", 179)
    ,("E160", 77)
    ,("A164", "Comoras")
    ,("B164", 197)
    ,("C164", 21)
    ,("D164", 127)
    ,("E164", 67)
    ,("G164", 1)
    ,("H164", 3)
    ,("A165", "Islas Caimanes")
    ,("D165", 115)
    ,("E165", 71)
    ,("H165", 1)
    ,("A166", "Islas Feroe")
    ,("B166", 187)
    ,("D166", 187)
    ,("E166", 0)
    ,("H166", 0)
    ,("A167", "Siria")
    ,("B167", 177)
    ,("D167", 78)
    ,("E167", 93)
    ,("H167", 6)
    ,("B176", 121)
    ,("C176", 12)
    ,("E176", 82)
    ,("A213", "Islas Virgenes Britanicas")
    ,("D213", 7)
    ,("H213", 1)
    ,("A214", "Papua Nueva Guinea")
    ,("D214", 8)
    ,("H214", 0)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
